$wb = $excel.ActiveWorkbook

# --- Second sheet: re-point selection to A2:A11 and apply British short-date
#     number format (dd/mm/yy) to the date column ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("A2:A11").NumberFormat = "dd/mm/yy;@"
$ws2.Range("A2:A11").Select()

# --- Rename "Sheet3" to "Sheet 3" ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Sheet 3"

# --- Sheet 3 becomes the active tab ---
$ws3.Activate()
